$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in cell values in the same order the author typed them, so that the
# shared-strings table gets built up in the matching order:
#   Row 82 -> Row 83 -> Row 84 -> Row 85 -> Row 87 -> Row 86

$ws.Range("A82").Value = "cut-and-dried"
$ws.Range("B82").Value = "قص و تجفيف"

$ws.Range("A83").Value = "reminiscing"
$ws.Range("B83").Value = "ذكريات"
$ws.Range("C83").Value = 33

$ws.Range("A84").Value = "mendous storage"
$ws.Range("B84").Value = "تخزين هائل"
$ws.Range("C84").Value = 33

$ws.Range("A85").Value = "consequently"
$ws.Range("B85").Value = "بناء على ذلك"
$ws.Range("C85").Value = 33

$ws.Range("A87").Value = "fraud detection"
$ws.Range("B87").Value = "الكشف عن الغش"
$ws.Range("C87").Value = 35

$ws.Range("A86").Value = "broken up"
$ws.Range("B86").Value = "إنفصلنا"
$ws.Range("C86").Value = 34

# Update the view to match the edited workbook state (scroll position +
# active cell selection, as Excel records when a user finishes editing here)
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 70
$win.ScrollColumn = 1
$ws.Range("A68").Select()
